$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Row 15 ---
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 58.333333333333
$ws.Range("L15").Value = 90
$ws.Range("M15").Value = 26.666666666666
$ws.Range("N15").Value = 72.727272727272

# --- Row 16 ---
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -28.571428571428
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 173
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = 57.272727272727
$ws.Range("L16").Value = 64.761904761904
$ws.Range("M16").Value = 0.581395348837
$ws.Range("N16").Value = -80.841638981173

# --- Row 17 ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 41.666666666666
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 184
$ws.Range("K17").Value = 8.695652173913
$ws.Range("L17").Value = 41.843971631205
$ws.Range("M17").Value = 70.94017094017
$ws.Range("N17").Value = -29.078014184397

# --- Row 18 ---
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -54.166666666666
$ws.Range("I18").Value = 147
$ws.Range("J18").Value = 169
$ws.Range("K18").Value = -13.017751479289
$ws.Range("L18").Value = -8.125
$ws.Range("M18").Value = -34.666666666666
$ws.Range("N18").Value = -89.33236574746

# --- Row 19 ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -6.521739130434
$ws.Range("I19").Value = 619
$ws.Range("J19").Value = 403
$ws.Range("K19").Value = 53.598014888337
$ws.Range("L19").Value = 65.066666666666
$ws.Range("M19").Value = 48.086124401913
$ws.Range("N19").Value = -25.15114873035

# --- Row 20 ---
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 13.333333333333
$ws.Range("I20").Value = 211
$ws.Range("J20").Value = 164
$ws.Range("K20").Value = 28.658536585365
$ws.Range("L20").Value = 35.25641025641
$ws.Range("M20").Value = 12.83422459893
$ws.Range("N20").Value = -88.788522848034

# --- Row 21 ---
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -25.714285714285
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -7.894736842105
$ws.Range("I21").Value = 1370
$ws.Range("J21").Value = 1044
$ws.Range("K21").Value = 31.226053639846
$ws.Range("L21").Value = 44.36248682824
$ws.Range("M21").Value = 20.598591549295
$ws.Range("N21").Value = -74.106974106974

# --- Row 22 ---
$ws.Range("F22").Value = 7
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 73
$ws.Range("K22").Value = 114.705882352941
$ws.Range("L22").Value = 135.483870967742
$ws.Range("M22").Value = 62.222222222222

# --- Row 24 ---
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 136
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 46.236559139784
$ws.Range("I24").Value = 1354
$ws.Range("J24").Value = 1136
$ws.Range("K24").Value = 19.19014084507
$ws.Range("L24").Value = 65.323565323565
$ws.Range("M24").Value = 60.616844602609

# --- Row 25 ---
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 11.111111111111
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 25.714285714285
$ws.Range("I25").Value = 505
$ws.Range("J25").Value = 425
$ws.Range("K25").Value = 18.823529411764
$ws.Range("L25").Value = 68.896321070234
$ws.Range("M25").Value = 12.723214285714

# --- Row 26 ---
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 25
$ws.Range("J26").Value = 21
$ws.Range("K26").Value = 19.047619047619
$ws.Range("L26").Value = 47.058823529411

# --- Row 27 ---
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 10
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 83
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = 45.614035087719
$ws.Range("L27").Value = 45.614035087719

# --- Row 28 ---
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 2
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = -50
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("F28").Value = 1
$ws.Range("G28").NumberFormat = '#,##0'
$ws.Range("G28").Value = 2
$ws.Range("H28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = 20
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = -57.142857142857

# --- Row 29 ---
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("C29").Value = 1
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E29").Value = 0
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 1
$ws.Range("G29").NumberFormat = '#,##0'
$ws.Range("G29").Value = 1
$ws.Range("H29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 5
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -28.571428571428
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 400
$ws.Range("N29").Value = -64.285714285714
